$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.195364238410596
$ws.Range("C2").Value = 0.5596026490066225
$ws.Range("J2").Value = 0.009933774834437087
$ws.Range("P2").Value = 0.1721854304635762
$ws.Range("S2").Value = 0.06291390728476821
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.0335195530726257
$ws.Range("J3").Value = 0.0335195530726257
$ws.Range("P3").Value = 0.7374301675977654
$ws.Range("S3").Value = 0.1899441340782123
$ws.Range("J4").Value = 0.02
$ws.Range("O4").Value = 0.02
$ws.Range("P4").Value = 0.66
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.06842105263157895
$ws.Range("D6").Value = 0.01052631578947368
$ws.Range("F6").Value = 0.05789473684210526
$ws.Range("J6").Value = 0.2578947368421053
$ws.Range("O6").Value = 0.01052631578947368
$ws.Range("Q6").Value = 0.1842105263157895
$ws.Range("R6").Value = 0.08421052631578947
$ws.Range("S6").Value = 0.3263157894736842
$ws.Range("B7").Value = 0.1294117647058824
$ws.Range("D7").Value = 0.01176470588235294
$ws.Range("E7").Value = 0.005882352941176471
$ws.Range("F7").Value = 0.04705882352941176
$ws.Range("J7").Value = 0.1529411764705882
$ws.Range("O7").Value = 0.02352941176470588
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.08235294117647059
$ws.Range("S7").Value = 0.3470588235294118
$ws.Range("B8").Value = 0.09212730318257957
$ws.Range("D8").Value = 0.01340033500837521
$ws.Range("E8").Value = 0.003350083752093802
$ws.Range("F8").Value = 0.04187604690117253
$ws.Range("J8").Value = 0.1222780569514238
$ws.Range("O8").Value = 0.01842546063651591
$ws.Range("Q8").Value = 0.202680067001675
$ws.Range("R8").Value = 0.08877721943048576
$ws.Range("S8").Value = 0.4170854271356784
$ws.Range("B9").Value = 0.1072961373390558
$ws.Range("D9").Value = 0.02575107296137339
$ws.Range("F9").Value = 0.05579399141630902
$ws.Range("J9").Value = 0.1158798283261803
$ws.Range("O9").Value = 0.01716738197424893
$ws.Range("Q9").Value = 0.2360515021459227
$ws.Range("R9").Value = 0.07725321888412018
$ws.Range("S9").Value = 0.3648068669527897
$ws.Range("B10").Value = 0.08344370860927153
$ws.Range("D10").Value = 0.02119205298013245
$ws.Range("E10").Value = 0.001986754966887417
$ws.Range("F10").Value = 0.05099337748344371
$ws.Range("J10").Value = 0.1298013245033112
$ws.Range("O10").Value = 0.01324503311258278
$ws.Range("Q10").Value = 0.2390728476821192
$ws.Range("R10").Value = 0.07748344370860927
$ws.Range("S10").Value = 0.3827814569536424
$ws.Range("G11").Value = 0.1575091575091575
$ws.Range("J11").Value = 0.1062271062271062
$ws.Range("K11").Value = 0.1978021978021978
$ws.Range("L11").Value = 0.5274725274725275
$ws.Range("S11").Value = 0.01098901098901099
$ws.Range("G12").Value = 0.7310344827586207
$ws.Range("J12").Value = 0.1724137931034483
$ws.Range("K12").Value = 0.02068965517241379
$ws.Range("L12").Value = 0.04137931034482759
$ws.Range("S12").Value = 0.03448275862068965
$ws.Range("F15").Value = 0.01219512195121951
$ws.Range("H15").Value = 0.1747967479674797
$ws.Range("I15").Value = 0.04471544715447155
$ws.Range("J15").Value = 0.3943089430894309
$ws.Range("K15").Value = 0.05691056910569105
$ws.Range("M15").Value = 0.01219512195121951
$ws.Range("O15").Value = 0.05691056910569105
$ws.Range("S15").Value = 0.2479674796747967
$ws.Range("F16").Value = 0.009345794392523364
$ws.Range("H16").Value = 0.1869158878504673
$ws.Range("I16").Value = 0.08878504672897196
$ws.Range("J16").Value = 0.4345794392523364
$ws.Range("K16").Value = 0.09345794392523364
$ws.Range("M16").Value = 0.01869158878504673
$ws.Range("O16").Value = 0.04205607476635514
$ws.Range("S16").Value = 0.1261682242990654
$ws.Range("F17").Value = 0.015
$ws.Range("H17").Value = 0.2416666666666667
$ws.Range("I17").Value = 0.08666666666666667
$ws.Range("J17").Value = 0.4033333333333333
$ws.Range("K17").Value = 0.06833333333333333
$ws.Range("M17").Value = 0.02166666666666667
$ws.Range("O17").Value = 0.05666666666666666
$ws.Range("S17").Value = 0.1066666666666667
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.212962962962963
$ws.Range("I18").Value = 0.06481481481481481
$ws.Range("J18").Value = 0.4398148148148148
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.03703703703703703
$ws.Range("O18").Value = 0.06018518518518518
$ws.Range("S18").Value = 0.06018518518518518
$ws.Range("F19").Value = 0.01
$ws.Range("H19").Value = 0.235
$ws.Range("I19").Value = 0.09857142857142857
$ws.Range("J19").Value = 0.3871428571428571
$ws.Range("K19").Value = 0.08071428571428571
$ws.Range("M19").Value = 0.01285714285714286
$ws.Range("N19").Value = 0.002142857142857143
$ws.Range("O19").Value = 0.07428571428571429
$ws.Range("S19").Value = 0.09928571428571428
